# Auto commit at 2025-11-13  8:05:27.77
#
# Appends the next day's (2025-11-12, Excel serial 45973) two station rows
# -- "四方坪站充电量(kw)" and "高岭站充电量(kw)" -- to the bottom of the daily
# charging-load log on Sheet1, mirroring the existing row 144/145 pattern,
# and moves the view/selection down to follow the newly entered data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dateSerial = 45973

$row146 = [ordered]@{
    A = $dateSerial
    B = "四方坪站充电量(kw)"
    C = 449.863
    D = 1237.5750000000003
    E = 392.76
    F = 269.76300000000003
    G = 343.61699999999996
    H = 587.90499999999997
    I = 358.90300000000002
    J = 158.786
    K = 148.13900000000001
    L = 237.64000000000001
    M = 291.92500000000001
    N = 175.49
    O = 982.87800000000016
    P = 1048.23
    Q = 304.76799999999997
    R = 443.279
    S = 209.01599999999999
    T = 146.59500000000003
    U = 75.387
    V = 173.78000000000003
    W = 69.33
    X = 144.94
    Y = 42.480000000000004
    Z = 72.89
}

$row147 = [ordered]@{
    A = $dateSerial
    B = "高岭站充电量(kw)"
    C = 210.30400000000003
    D = 402.40699999999998
    E = 138.82
    F = 143.86799999999999
    G = 0
    H = 102.57499999999999
    I = 354.30099999999999
    J = 265.68099999999998
    K = 172.53800000000001
    L = 94.256
    M = 239.815
    N = 240.05500000000001
    O = 400.279
    P = 418.98099999999999
    Q = 372.07800000000003
    R = 224.238
    S = 221.90499999999997
    T = 199.76100000000002
    U = 116.20400000000001
    V = 56.335999999999999
    W = 98.685999999999993
    X = 39.241999999999997
    Y = 50.230000000000004
    Z = 7.0430000000000001
}

foreach ($col in $row146.Keys) {
    $ws.Range($col + "146").Value = $row146[$col]
}
foreach ($col in $row147.Keys) {
    $ws.Range($col + "147").Value = $row147[$col]
}

# Keep the workbook scrolled/selected near the freshly appended rows, same
# as the source workbook after a new day's rows are pasted in.
$ws.Range("A134").Select() | Out-Null
$ws.Range("B149").Select() | Out-Null
